# Applies the coin price / volume(1h) refresh described in the commit
# message ("Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated Price (column D) / Volume(1h) (column E) text.
# A row is omitted from a key when that column did not change for that row.
$rowUpdates = @(
    @{ Row = 2; D = "28.484.87"; E = "  -0.21%  " },
    @{ Row = 3; D = "1.868.98"; E = "  -0.44%  " },
    @{ Row = 4; D = "1.010"; E = "  -1.37%  " },
    @{ Row = 5; D = "315.28"; E = "  -1.03%  " },
    @{ Row = 6; E = "  -1.48%  " },
    @{ Row = 7; D = "0.5077"; E = "  -1.50%  " },
    @{ Row = 8; D = "0.3907"; E = "  -1.62%  " },
    @{ Row = 9; D = "0.08353"; E = "  -0.23%  " },
    @{ Row = 10; D = "42.45"; E = "  +0.62%  " },
    @{ Row = 11; D = "1.104"; E = "  -0.93%  " },
    @{ Row = 12; D = "6.202"; E = "  -1.01%  " },
    @{ Row = 13; D = "1.867.27"; E = "  +2.58%  " },
    @{ Row = 14; D = "20.33"; E = "  -1.14%  " },
    @{ Row = 15; D = "7.242"; E = "  -0.04%  " },
    @{ Row = 16; D = "1.010"; E = "  -1.25%  " },
    @{ Row = 17; D = "0.00001100"; E = "  -1.03%  " },
    @{ Row = 18; D = "91.28"; E = "  -0.12%  " },
    @{ Row = 19; D = "0.06724"; E = "  -0.87%  " },
    @{ Row = 20; D = "17.64"; E = "  -0.74%  " },
    @{ Row = 21; D = "1.008"; E = "  -1.48%  " },
    @{ Row = 22; D = "5.903"; E = "  -1.24%  " },
    @{ Row = 23; D = "28.548.22"; E = "  -0.11%  " },
    @{ Row = 24; E = "  -0.91%  " },
    @{ Row = 25; D = "2.203"; E = "  -3.92%  " },
    @{ Row = 26; D = "2.082.03"; E = "  +2.34%  " },
    @{ Row = 27; D = "157.06"; E = "  -3.38%  " },
    @{ Row = 28; D = "20.55"; E = "  -1.01%  " },
    @{ Row = 29; D = "2.420"; E = "  +2.36%  " },
    @{ Row = 30; D = "125.93"; E = "  -1.47%  " },
    @{ Row = 31; E = "  -1.48%  " },
    @{ Row = 32; D = "1.038"; E = "  -0.39%  " },
    @{ Row = 33; D = "5.752"; E = "  -1.37%  " },
    @{ Row = 34; D = "3.619"; E = "  -0.69%  " },
    @{ Row = 35; D = "0.02456"; E = "  +0.79%  " },
    @{ Row = 36; D = "0.06604"; E = "  +1.41%  " },
    @{ Row = 37; D = "9.019"; E = "  +1.15%  " },
    @{ Row = 38; D = "0.2159"; E = "  -1.49%  " },
    @{ Row = 39; D = "5.036" },
    @{ Row = 40; D = "1.180"; E = "  -0.84%  " },
    @{ Row = 41; D = "1.236"; E = "  -3.30%  " },
    @{ Row = 42; D = "0.6359"; E = "  -1.37%  " },
    @{ Row = 43; D = "11.08"; E = "  -1.73%  " },
    @{ Row = 44; E = "  -1.40%  " },
    @{ Row = 45; D = "0.5995"; E = "  -0.83%  " },
    @{ Row = 46; D = "13.08"; E = "  +0.32%  " },
    @{ Row = 47; D = "3.680"; E = "  -1.50%  " },
    @{ Row = 48; D = "1.998"; E = "  -0.01%  " },
    @{ Row = 49; D = "1.211"; E = "  +0.06%  " },
    @{ Row = 50; D = "122.36"; E = "  +0.18%  " },
    @{ Row = 51; D = "1.122" }
)

foreach ($u in $rowUpdates) {
    if ($u.ContainsKey("D")) {
        # Price values such as "1.010" or "0.00001100" look numeric to Excel
        # and would otherwise be parsed as a Double (losing the exact text,
        # e.g. trailing zeros, or thousands-grouped values like "28.484.87").
        # Force the cell to Text format, assign the literal string, then put
        # the style back so no stray formatting is left behind.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        # Volume(1h) values are padded with spaces and a trailing "%", so
        # Excel always stores them as plain text already.
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

